$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cell values from the crypto price refresh.
# NumberFormat "@" + ClearFormats() ensures digit-like strings (e.g. "573.83")
# are stored as text, matching the original inlineStr cell type, without leaving
# a residual style index behind.

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '63.834.21'
$cell.ClearFormats()

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '2.538.95'
$cell.ClearFormats()

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  +5.68%  '
$cell.ClearFormats()

$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.ClearFormats()

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -0.06%  '
$cell.ClearFormats()

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '573.83'
$cell.ClearFormats()

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  +2.38%  '
$cell.ClearFormats()

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '146.78'
$cell.ClearFormats()

$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  +5.98%  '
$cell.ClearFormats()

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.ClearFormats()

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -0.05%  '
$cell.ClearFormats()

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  +0.50%  '
$cell.ClearFormats()

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '2.541.41'
$cell.ClearFormats()

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  +5.89%  '
$cell.ClearFormats()

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  +2.21%  '
$cell.ClearFormats()

$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '5.76'
$cell.ClearFormats()

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  +0.85%  '
$cell.ClearFormats()

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  +1.70%  '
$cell.ClearFormats()

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '0.358'
$cell.ClearFormats()

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  +2.95%  '
$cell.ClearFormats()

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '27.76'
$cell.ClearFormats()

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  +8.29%  '
$cell.ClearFormats()

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '2.991.09'
$cell.ClearFormats()

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  +5.67%  '
$cell.ClearFormats()

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '63.613.10'
$cell.ClearFormats()

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +2.83%  '
$cell.ClearFormats()

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  +4.07%  '
$cell.ClearFormats()

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '2.535.68'
$cell.ClearFormats()

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  +5.09%  '
$cell.ClearFormats()

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  +4.28%  '
$cell.ClearFormats()

$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '342.96'
$cell.ClearFormats()

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  +0.42%  '
$cell.ClearFormats()

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  +3.11%  '
$cell.ClearFormats()

$cell = $ws.Range('E22')
$cell.NumberFormat = '@'
$cell.Value = '  +0.11%  '
$cell.ClearFormats()

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  +0.59%  '
$cell.ClearFormats()

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '65.98'
$cell.ClearFormats()

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  +1.79%  '
$cell.ClearFormats()

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -0.37%  '
$cell.ClearFormats()

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '1.58'
$cell.ClearFormats()

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  +5.98%  '
$cell.ClearFormats()

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  +0.07%  '
$cell.ClearFormats()

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '8.25'
$cell.ClearFormats()

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  -0.85%  '
$cell.ClearFormats()

$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '1.43'
$cell.ClearFormats()

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  +4.01%  '
$cell.ClearFormats()

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0825'
$cell.ClearFormats()

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  +7.77%  '
$cell.ClearFormats()

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  +4.19%  '
$cell.ClearFormats()

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '6.83'
$cell.ClearFormats()

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  +7.53%  '
$cell.ClearFormats()

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '176.99'
$cell.ClearFormats()

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  +3.30%  '
$cell.ClearFormats()

$cell = $ws.Range('B34')
$cell.NumberFormat = '@'
$cell.Value = 'Bittensor'
$cell.ClearFormats()

$cell = $ws.Range('C34')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell.ClearFormats()

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '417.05'
$cell.ClearFormats()

$cell = $ws.Range('E34')
$cell.NumberFormat = '@'
$cell.Value = '  +14.66%  '
$cell.ClearFormats()

$cell = $ws.Range('B35')
$cell.NumberFormat = '@'
$cell.Value = 'ImmutableX'
$cell.ClearFormats()

$cell = $ws.Range('C35')
$cell.NumberFormat = '@'
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.ClearFormats()

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '1.54'
$cell.ClearFormats()

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  +10.32%  '
$cell.ClearFormats()

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.401'
$cell.ClearFormats()

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  +2.16%  '
$cell.ClearFormats()

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '19.07'
$cell.ClearFormats()

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  +3.33%  '
$cell.ClearFormats()

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  -2.73%  '
$cell.ClearFormats()

$cell = $ws.Range('E39')
$cell.NumberFormat = '@'
$cell.Value = '  +0.01%  '
$cell.ClearFormats()

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  +6.11%  '
$cell.ClearFormats()

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.998'
$cell.ClearFormats()

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  -0.14%  '
$cell.ClearFormats()

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '40.71'
$cell.ClearFormats()

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  +4.59%  '
$cell.ClearFormats()

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '152.74'
$cell.ClearFormats()

$cell = $ws.Range('E43')
$cell.NumberFormat = '@'
$cell.Value = '  +6.22%  '
$cell.ClearFormats()

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '3.78'
$cell.ClearFormats()

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  +3.59%  '
$cell.ClearFormats()

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '20.90'
$cell.ClearFormats()

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  +2.66%  '
$cell.ClearFormats()

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  +2.70%  '
$cell.ClearFormats()

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '0.0966'
$cell.ClearFormats()

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  +0.76%  '
$cell.ClearFormats()

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '19.03'
$cell.ClearFormats()

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  +7.28%  '
$cell.ClearFormats()

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0232'
$cell.ClearFormats()

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  +4.66%  '
$cell.ClearFormats()

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  +7.39%  '
$cell.ClearFormats()
